$d = $word.ActiveDocument

# --- Change 1: merge the split "graf" / "ic" runs into a single "grafic" run ---
# (inside "Descriere toolbox grafic Labview"; this is the only place in the
# document where "graf" and "ic" appear as adjacent separate runs/words)
$d.Content.Find.Execute("graf" + "ic", $true, $false, $false, $false, $false, `
    $true, 1, $false, "grafic", 2) | Out-Null

# --- Change 2: fix "Toolbox graphic Labview" -> "Toolbox grafic Labview" and
#     relocate the singleton _GoBack bookmark into the middle of that word ---
$titlePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.TrimEnd([char]13) -eq "Toolbox graphic Labview") {
        $titlePara = $cand
    }
}
$tStart = $titlePara.Range.Start

# "Toolbox graphic " offsets relative to the paragraph start:
#   0-10  "Toolbox gra"
#   11-12 "ph"   -> becomes "f"
#   13-15 "ic "
#   16-22 "Labview"
#
# Step 2a: split the run that carries the lastRenderedPageBreak marker away
# from the rest of the word by placing the (singleton, auto-relocating)
# _GoBack bookmark right after "Toolbox gra". Because this only inserts a
# bookmark (no text is rewritten) the original run - and its page-break
# marker - survive untouched.
$splitPoint = $d.Range($tStart + 11, $tStart + 11)
$d.Bookmarks.Add("_GoBack", $splitPoint) | Out-Null

# Step 2b: turn "ph" into "f". This now only touches the freshly split-off
# run (which holds no special markers), so nothing else is disturbed.
$phRange = $d.Range($tStart + 11, $tStart + 13)
$phRange.Text = "f"

# Step 2c: move the _GoBack bookmark to its real target spot, between "f"
# and "ic ". Re-adding the singleton bookmark relocates it automatically and
# splits the "fic " run into "f" | "ic " in the process.
$goBackPoint = $d.Range($tStart + 12, $tStart + 12)
$d.Bookmarks.Add("_GoBack", $goBackPoint) | Out-Null

# --- Change 3: drop the stray "sda" paragraph (its old _GoBack bookmark has
#     already been relocated above, so nothing is orphaned) ---
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13) -eq "sda") {
        $p.Range.Delete()
    }
}
